$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update account IDs in column A (new numbering scheme HD1130000X)
$ws.Range("A1").Value = "HD11300001"
$ws.Range("A2").Value = "HD11300002"
$ws.Range("A3").Value = "HD11300003"
$ws.Range("A4").Value = "HD11300004"
$ws.Range("A5").Value = "HD11300005"
$ws.Range("A6").Value = "HD11300006"
$ws.Range("A7").Value = "HD11300007"

# Update passwords for admin accounts (validate data when login and change password)
$ws.Range("B8").Value = "password123"
$ws.Range("B9").Value = "ngochoai123"
$ws.Range("B10").Value = "bichdiep123"

# Widen column B to fit the new longer passwords
$ws.Columns.Item(2).ColumnWidth = 18.83

# Move the active selection to E11
$ws.Range("E11").Select()
